$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.364.05'
$ws.Range("E2").Value = '  -0.85%  '
$ws.Range("D3").Value = '2.178.07'
$ws.Range("E3").Value = '  -2.06%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.52%  '
$ws.Range("E6").Value = '  -1.31%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '72.63'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.21%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.580'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.92'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0908'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.43%  '
$ws.Range("E12").Value = '  -0.23%  '
$ws.Range("E13").Value = '  -2.65%  '
$ws.Range("D14").Value = '2.509.53'
$ws.Range("E14").Value = '  -1.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.10'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.19%  '
$ws.Range("D16").Value = '2.175.33'
$ws.Range("E16").Value = '  -1.81%  '
$ws.Range("E17").Value = '  -4.54%  '
$ws.Range("D18").Value = '42.298.76'
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000101'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.62%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.54'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.45%  '
$ws.Range("E21").Value = '  -2.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.12'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.72%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.31'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.80%  '
$ws.Range("E24").Value = '  -2.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.40'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.36'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.21'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.81%  '
$ws.Range("E29").Value = '  -2.66%  '
$ws.Range("E30").Value = '  -1.45%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '36.34'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.96'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0811'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.07'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.63%  '
$ws.Range("E35").Value = '  -1.76%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.106'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.17'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0335'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.66'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.24%  '
$ws.Range("E40").Value = '  -4.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.194'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.68%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '59.03'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.11'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.06'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.45'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.95%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0970'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.40%  '
$ws.Range("B47").Value = 'WOONetwork'
$ws.Range("C47").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.457'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.74%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.16'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.07'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.23%  '
$ws.Range("E50").Value = '  -1.41%  '
$ws.Range("E51").Value = '  +0.24%  '
